$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a numeric-looking string to be stored as exact text (avoid Excel coercing it to a Double,
# which would lose trailing zeros / introduce floating point noise), while keeping the default cell style.
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '65.855.05'
$ws.Range('E2').Value = '  -1.86%  '
$ws.Range('D3').Value = '3.437.76'
$ws.Range('E3').Value = '  -0.51%  '
Set-TextValue 'D5' '583.12'
$ws.Range('E5').Value = '  -0.28%  '
Set-TextValue 'D6' '173.32'
$ws.Range('E6').Value = '  -1.60%  '
$ws.Range('E7').Value = '  -0.03%  '
Set-TextValue 'D8' '0.604'
$ws.Range('E8').Value = '  -0.72%  '
$ws.Range('D9').Value = '3.434.95'
$ws.Range('E9').Value = '  -0.57%  '
$ws.Range('E10').Value = '  -2.99%  '
Set-TextValue 'D11' '6.93'
$ws.Range('E11').Value = '  -0.12%  '
Set-TextValue 'D12' '0.409'
$ws.Range('E12').Value = '  -3.49%  '
$ws.Range('D13').Value = '4.030.58'
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('E14').Value = '  +1.44%  '
Set-TextValue 'D15' '28.68'
$ws.Range('E15').Value = '  -9.27%  '
$ws.Range('D16').Value = '65.908.80'
$ws.Range('E17').Value = '  -2.15%  '
$ws.Range('D18').Value = '3.438.99'
$ws.Range('E18').Value = '  -0.83%  '
Set-TextValue 'D19' '5.93'
$ws.Range('E19').Value = '  -2.18%  '
$ws.Range('E20').Value = '  -0.48%  '
Set-TextValue 'D21' '368.65'
$ws.Range('E21').Value = '  -2.82%  '
$ws.Range('E22').Value = '  -1.96%  '
Set-TextValue 'D23' '72.15'
$ws.Range('E23').Value = '  +0.92%  '
Set-TextValue 'D24' '0.999'
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('E26').Value = '  +0.82%  '
$ws.Range('E27').Value = '  -1.95%  '
Set-TextValue 'D28' '0.177'
$ws.Range('E28').Value = '  +1.67%  '
Set-TextValue 'D29' '0.999'
$ws.Range('E29').Value = '  -0.18%  '
Set-TextValue 'D30' '23.63'
$ws.Range('E30').Value = '  -1.31%  '
$ws.Range('E31').Value = '  -3.98%  '
Set-TextValue 'D32' '1.98'
$ws.Range('E32').Value = '  -2.20%  '
Set-TextValue 'D33' '0.999'
$ws.Range('E33').Value = '  -0.01%  '
$ws.Range('E34').Value = '  -5.71%  '
Set-TextValue 'D35' '7.00'
$ws.Range('E35').Value = '  -2.67%  '
$ws.Range('E36').Value = '  -0.45%  '
Set-TextValue 'D37' '160.55'
$ws.Range('E37').Value = '  +0.28%  '
$ws.Range('B38').Value = 'EnergySwap'
$ws.Range('C38').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D38' '28.70'
$ws.Range('E38').Value = '  +5.88%  '
$ws.Range('B39').Value = 'Mantle'
$ws.Range('C39').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D39' '0.878'
$ws.Range('E39').Value = '  +0.16%  '
$ws.Range('E40').Value = '  -2.84%  '
Set-TextValue 'D41' '2.60'
$ws.Range('E41').Value = '  -1.54%  '
$ws.Range('D42').Value = '2.772.29'
$ws.Range('E42').Value = '  +2.87%  '
$ws.Range('E43').Value = '  -2.72%  '
$ws.Range('E44').Value = '  -0.42%  '
Set-TextValue 'D45' '0.0681'
$ws.Range('E45').Value = '  -2.32%  '
Set-TextValue 'D46' '40.13'
$ws.Range('E46').Value = '  -2.27%  '
Set-TextValue 'D47' '24.48'
$ws.Range('E47').Value = '  -3.59%  '
$ws.Range('E48').Value = '  -1.63%  '
Set-TextValue 'D49' '324.36'
$ws.Range('E49').Value = '  +0.71%  '
$ws.Range('E50').Value = '  -1.42%  '
Set-TextValue 'D51' '6.25'
$ws.Range('E51').Value = '  +0.36%  '
